# Commit: "it should create a new group when payer login with a single clone scene"
#
# The Scene config sheet has a "CanClone" column (K). For the single-clone
# scenes (rows 10, 11, 13, 14, 15) CanClone was incorrectly set to 1 (true).
# Flip it to 0 (false) so a payer logging into one of these single-clone
# scenes creates a brand-new group instead of being allowed to clone/share
# an existing one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K10").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("K15").Value = 0

# Leave the active selection on the first edited cell, matching the
# author's saved view state.
$ws.Range("K10").Select()
